$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy formatting from the row above (A53) into the new row's date cell (A54),
# then set the new values for the appended row 54.
$ws.Range("A53").Copy($ws.Range("A54"))
$ws.Range("A54").Value = 45986

$ws.Range("B54").Value = 2025
$ws.Range("C54").Value = -0.08656168856399082
$ws.Range("D54").Value = 2026
$ws.Range("E54").Value = 0.6232357314897463
